# Adds a new "2022-Q3" sheet (with fresh fund-holder data) right after the
# "总计" summary sheet / right before "2022-Q2", and records the new
# quarter's totals in the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: force a cell's value to be stored as TEXT (not auto-coerced to
# a number), mirroring how the source file keeps numeric-looking strings
# like "0.55" / "010783" as inline/shared strings rather than numbers.
# Setting NumberFormat="@" first forces text interpretation; pasting the
# *format only* from a plain, already-default-styled cell afterwards
# wipes out the incidental style (quote-prefix / text-format) that step
# introduces, so the cell ends up looking exactly like its neighbours.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($Cell, $Text, $FormatDonor)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    if ($FormatDonor) {
        $FormatDonor.Copy() | Out-Null
        $Cell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    }
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by cloning "2022-Q2" (same column
#    layout / header / 6-row shape / cell styles), then overwrite the
#    data with the new quarter's numbers.
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)                       # places the clone right before 2022-Q2
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

$donor = $wsQ3.Range("C2")              # a plain, default-styled text cell

$q3data = @(
    @("010783", "德邦沪港深龙头混合A",                         "0.55", "84.96", "3.08", "0.0169", 9),
    @("013897", "德邦港股通成长精选混合型证券投资基金A",         "0.41", "79.99", "3.05", "0.0125", 9),
    @("013898", "德邦港股通成长精选混合型证券投资基金C",         "0.37", "79.99", "3.05", "0.0113", 9),
    @("010784", "德邦沪港深龙头混合C",                         "0.36", "84.96", "3.08", "0.0111", 9),
    @("005143", "中融沪港深大消费主题灵活配置混合C",             "0.27", "90.10", "3.94", "0.0106", 10),
    @("005142", "中融沪港深大消费主题灵活配置混合A",             "0.13", "90.10", "3.94", "0.0051", 10)
)

for ($i = 0; $i -lt $q3data.Length; $i++) {
    $row = $i + 2
    $vals = $q3data[$i]

    Set-TextValue $wsQ3.Range("B$row") $vals[0] $donor
    $wsQ3.Range("C$row").Value = $vals[1]
    Set-TextValue $wsQ3.Range("D$row") $vals[2] $donor
    Set-TextValue $wsQ3.Range("E$row") $vals[3] $donor
    Set-TextValue $wsQ3.Range("F$row") $vals[4] $donor
    Set-TextValue $wsQ3.Range("G$row") $vals[5] $donor
    $wsQ3.Range("H$row").Value = $vals[6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert the 2022-Q3 row at the top
#    of the data, shift every following quarter down one row, and add
#    the row for the quarter that now falls off the end (2020-Q4).
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$totalData = @(
    @("2022-Q3", 6,  0.07000000000000001),
    @("2022-Q2", 6,  0.89),
    @("2022-Q1", 5,  0.91),
    @("2021-Q4", 9,  3.82),
    @("2021-Q3", 20, 8.68),
    @("2021-Q2", 25, 11.93),
    @("2021-Q1", 8,  3.72),
    @("2020-Q4", 5,  1.82)
)

# Column-A style donor (bold/centered/bordered "index" style already used
# by the existing rows) so the brand-new last row matches its siblings.
$aDonor = $wsTotal.Range("A2")

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $row = $i + 2
    $vals = $totalData[$i]

    if ($row -eq ($totalData.Length + 1)) {
        # shouldn't happen, guard kept for clarity
    }

    $aCell = $wsTotal.Range("A$row")
    if ($row -gt 8) {
        $aDonor.Copy() | Out-Null
        $aCell.PasteSpecial(-4122) | Out-Null
    }
    $aCell.Value = $i

    $wsTotal.Range("B$row").Value = $vals[0]
    $wsTotal.Range("C$row").Value = $vals[1]
    $wsTotal.Range("D$row").Value = $vals[2]
}

Write-Output "done"
